$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.755599999999998
$ws.Range("B3").Value = 6.005000000000003
$ws.Range("E3").Value = 15.86019999999999
$ws.Range("E12").Value = 17.19310000000003
$ws.Range("B14").Value = 5.2077
$ws.Range("B21").Value = 9.365000000000004
$ws.Range("B23").Value = 9.096199999999998
$ws.Range("E24").Value = 16.7485
$ws.Range("B25").Value = 5.193700000000002
$ws.Range("D25").Value = -7.269199999999995
$ws.Range("E25").Value = 17.32210000000002
$ws.Range("B26").Value = 5.159600000000007
$ws.Range("D27").Value = -8.765300000000005
$ws.Range("B29").Value = 5.007500000000004
$ws.Range("D31").Value = -8.751200000000004
$ws.Range("D39").Value = -7.979700000000004
$ws.Range("D48").Value = -7.300399999999994
$ws.Range("E50").Value = 16.2934
$ws.Range("D51").Value = -7.376699999999996
$ws.Range("D52").Value = -7.616099999999998
$ws.Range("B53").Value = 5.150499999999998
$ws.Range("E53").Value = 16.79450000000001
$ws.Range("D55").Value = -8.415799999999997
$ws.Range("D56").Value = -7.6626
$ws.Range("B57").Value = 5.067899999999995
$ws.Range("D57").Value = -7.738500000000001
$ws.Range("E57").Value = 16.7438
$ws.Range("B59").Value = 4.895599999999996
$ws.Range("E61").Value = 16.55
$ws.Range("E63").Value = 17.45020000000003
$ws.Range("B69").Value = 5.308499999999998
$ws.Range("E70").Value = 17.28220000000003
$ws.Range("D73").Value = -7.925599999999996
$ws.Range("B79").Value = 9.066600000000005
$ws.Range("B83").Value = 5.189699999999999
$ws.Range("E86").Value = 16.5594
$ws.Range("D89").Value = -5.654900000000001
$ws.Range("D90").Value = -8.2475
$ws.Range("B91").Value = 4.919899999999998
$ws.Range("D92").Value = -5.732400000000001
$ws.Range("B93").Value = 5.744499999999999
$ws.Range("E98").Value = 15.65120000000001
$ws.Range("E100").Value = 16.5893
$ws.Range("E102").Value = 16.65749999999998
